# "Generate Report for Handoff"
# The localization status report moves from "In Translation" to
# "Ready for handoff" and the handoff timestamps are refreshed.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Status: "In Translation" -> "Ready for handoff" everywhere it appears.
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date / Latest Handoff Datetime: bump timestamps.
$ws1.Range("G2").Value = "2016-08-27 22:39:37"
$ws3.Range("H2").Value = "2016-08-27 22:39:37"
$ws2.Range("H2").Value = "2016-08-27 22:39:33"

# The wider "Ready for handoff" text makes the Status columns grow.
$ws1.Columns.Item(5).ColumnWidth = 16.3333333333333
$ws1.Columns.Item(6).ColumnWidth = 16.3333333333333
$ws2.Columns.Item(3).ColumnWidth = 16.3333333333333
$ws3.Columns.Item(3).ColumnWidth = 16.3333333333333
